# Auto-generated Excel COM-interop script to update cryptos list
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "26.144.37"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -0.01%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.654.68"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -0.14%  "

# Row 4
$ws.Range("E4").Value = "  -0.17%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "218.76"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.01%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.5240"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.17%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06353"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.94%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "20.53"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.30%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07691"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -1.43%  "

# Row 12
$ws.Range("B12").Value = "Polkadot"
$ws.Range("C12").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "4.619"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +2.74%  "

# Row 13
$ws.Range("B13").Value = "WrappedEther"
$ws.Range("C13").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.644.32"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.66%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "1.883.71"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.03%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.5613"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +1.01%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.0₅8199"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +2.36%  "

# Row 17
$ws.Range("E17").Value = "  +0.62%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "26.141.03"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.05%  "

# Row 19
$ws.Range("E19").Value = "  -0.20%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "4.655"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.37%  "

# Row 21
$ws.Range("E21").Value = "  +4.10%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "192.27"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -1.55%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.956"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.21%  "

# Row 24
$ws.Range("E24").Value = "  -0.21%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "144.90"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -1.23%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.1195"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.95%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "7.265"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +1.33%  "

# Row 28
$ws.Range("E28").Value = "  +0.01%  "

# Row 29
$ws.Range("E29").Value = "  +0.71%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.05448"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -5.00%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.272"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.13%  "

# Row 32
$ws.Range("E32").Value = "  -0.58%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.369"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.68%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.565"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -1.48%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.9535"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.02%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.779"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.84%  "

# Row 37
$ws.Range("E37").Value = "  -0.63%  "

# Row 38
$ws.Range("E38").Value = "  -0.53%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01585"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.76%  "

# Row 40
$ws.Range("E40").Value = "  -1.19%  "

# Row 41
$ws.Range("E41").Value = "  -0.16%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.8327"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -1.55%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.027.82"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -3.23%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "101.42"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -1.90%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.794.23"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.01%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "57.72"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.29%  "

# Row 47
$ws.Range("B47").Value = "Frax"
$ws.Range("C47").Value = "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.002"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.68%  "

# Row 48
$ws.Range("B48").Value = "EnergySwap"
$ws.Range("C48").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "8.001"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.17%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.4340"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -1.37%  "

# Row 50
$ws.Range("B50").Value = "Cronos"
$ws.Range("C50").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.05194"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -3.40%  "

# Row 51
$ws.Range("B51").Value = "BabyDogeCoin"
$ws.Range("C51").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0₈101"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -5.62%  "
